$wb = $excel.ActiveWorkbook

# --- validLogin sheet: selection moves from C3 to B1, tab no longer selected ---
$ws1 = $wb.Worksheets.Item("validLogin")
$ws1.Range("B1").Select()

# --- addnewcustomer sheet: new customer record, becomes the active tab ---
$ws3 = $wb.Worksheets.Item("addnewcustomer")
$ws3.Activate()
$ws3.Range("A1").Select()

# Update username / email for the two rows (new customer values)
$ws3.Range("H1").Value = "manoj24@gmail.com"
$ws3.Range("B1").Copy()
$ws3.Range("H1").PasteSpecial(-4122)

$ws3.Range("H2").Value = "manoj25@gmail.com"
$ws3.Range("B2").Copy()
$ws3.Range("H2").PasteSpecial(-4122)

$ws3.Range("A1").Value = "manojkumari"
$ws3.Range("A2").Value = "manojkumarj"

# Widen column H to fit the longer e-mail addresses
$ws3.Columns.Item(8).ColumnWidth = 19.6640625

$excel.CutCopyMode = $false
